$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.039475321769714
$ws.Range("B1").Value = 1.719291210174561
$ws.Range("C1").Value = 4.475198745727539
$ws.Range("D1").Value = 2.480980396270752
$ws.Range("E1").Value = 1.313926219940186
